# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / percentage / coin-name / URL cells (safe as direct string assignment) ---
$ws.Range("D2").Value = "65.479.36"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "3.396.37"
$ws.Range("E3").Value = "  -2.44%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("E6").Value = "  -4.22%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.397.29"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("E9").Value = "  -2.88%  "
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("E11").Value = "  -6.48%  "
$ws.Range("E12").Value = "  -4.99%  "
$ws.Range("D13").Value = "3.971.72"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("E14").Value = "  -7.06%  "
$ws.Range("E15").Value = "  -7.08%  "
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "65.546.95"
$ws.Range("D18").Value = "3.396.47"
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("E20").Value = "  -5.95%  "
$ws.Range("E21").Value = "  -5.82%  "
$ws.Range("E22").Value = "  -5.98%  "
$ws.Range("E23").Value = "  -5.53%  "
$ws.Range("E24").Value = "  -3.36%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "3.535.50"
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("E27").Value = "  -10.01%  "
$ws.Range("E28").Value = "  -6.17%  "
$ws.Range("E29").Value = "  -7.62%  "
$ws.Range("E30").Value = "  -3.29%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  -5.60%  "
$ws.Range("E33").Value = "  -8.77%  "
$ws.Range("E34").Value = "  -4.21%  "
$ws.Range("D35").Value = "3.392.30"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -9.21%  "
$ws.Range("E38").Value = "  -7.48%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E40").Value = "  -5.72%  "
$ws.Range("E41").Value = "  -4.47%  "
$ws.Range("E42").Value = "  -4.62%  "
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("E44").Value = "  -7.70%  "
$ws.Range("E45").Value = "  -10.41%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("E47").Value = "  -9.23%  "
$ws.Range("E48").Value = "  -4.91%  "
$ws.Range("E49").Value = "  -6.20%  "
$ws.Range("E50").Value = "  -7.59%  "
$ws.Range("E51").Value = "  -7.09%  "

# --- Numeric-looking price cells: force text formatting first so trailing zeros / multi-dot
#     "thousand-dot" formatted numbers are preserved exactly as text, then restore default style ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "413.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.160"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "168.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.868"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.01"
$ws.Range("D49").Style = "Normal"
